# "Trader changed to Edge" - update the robot config values that pointed at
# the old Trader mailbox / EARTH.GSI shared-drive paths so they point at the
# new mailbox and the new local ("Edge" machine) folder layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# TraderMailbox (row 33) - shared mailbox changed from the generic Trader
# mailbox to Laura Battell's mailbox.
$ws.Range("B33").Value = "Laura.Battell@defra.gov.uk"

# ToBePrintedFolder (row 29)
$ws.Range("B29").Value = "C:\Users\x953922\Desktop\CoFS\To Be Printed\"

# CertificateTemplatesPath (row 28)
$ws.Range("B28").Value = "C:\Users\x953922\Desktop\CoFS\Robot Certificate Templates\"

# EmailApplicationsFolder (row 27)
$ws.Range("B27").Value = "C:\Users\x953922\Desktop\CoFS\"

# Reflect the author's last on-screen selection/scroll position when they
# saved the workbook (row 27, cell B27).
$ws.Activate()
$ws.Range("B27").Select()
